# Rotates the three species-observation rows (2, 3, 4):
#   new row 2 <- old row 3
#   new row 3 <- old row 4
#   new row 4 <- old row 2
# Implemented as direct cell writes (rather than a physical row move) so
# that only the cells whose value actually changes are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cell updates ---
$numericUpdates = @(
  @{addr="A2"; val=94499215},
  @{addr="B2"; val=56859},
  @{addr="E2"; val=103018},
  @{addr="Q2"; val=722412.1116719937},
  @{addr="R2"; val=6640739.120328019},
  @{addr="S2"; val=10},

  @{addr="A3"; val=95907760},
  @{addr="B3"; val=55649},
  @{addr="E3"; val=208255},

  @{addr="A4"; val=72871343},
  @{addr="B4"; val=99590},
  @{addr="E4"; val=221333},
  @{addr="Q4"; val=721706.6784319634},
  @{addr="R4"; val=6640148.217632387},
  @{addr="S4"; val=5}
)
foreach ($u in $numericUpdates) {
    $ws.Range($u.addr).Value = $u.val
}

# --- Text cell updates ---
# NumberFormat is forced to "@" (text) before the assignment so values that
# look like numbers/dates ("1", "2021-06-26", ...) are not silently
# reinterpreted by Excel, then the format is cleared again so the cell is
# left without any explicit style, matching the rest of the sheet.
$textUpdates = @(
  @{addr='F2'; val='Svartvit flugsnappare'},
  @{addr='G2'; val='Ficedula hypoleuca'},
  @{addr='H2'; val='(Pallas, 1764)'},
  @{addr='I2'; val='1'},
  @{addr='Y2'; val='2021-06-26'},
  @{addr='AA2'; val='2021-06-26'},
  @{addr='AW2'; val='Jacob Törngren'},
  @{addr='AX2'; val='Jacob Törngren'},

  @{addr='D3'; val='LC'},
  @{addr='F3'; val='Skogsödla'},
  @{addr='G3'; val='Zootoca vivipara'},
  @{addr='H3'; val='(Jacquin, 1787)'},
  @{addr='I3'; val='2'},
  @{addr='J3'; val='ex.'},
  @{addr='Y3'; val='2021-09-04'},
  @{addr='AA3'; val='2021-09-04'},

  @{addr='D4'; val='NT'},
  @{addr='F4'; val='Backklöver'},
  @{addr='G4'; val='Trifolium montanum'},
  @{addr='H4'; val='L.'},
  @{addr='Y4'; val='2018-06-13'},
  @{addr='AA4'; val='2018-06-13'},
  @{addr='AH4'; val='Vägkant'},
  @{addr='AW4'; val='Johan Lilja'},
  @{addr='AX4'; val='Johan Lilja'},
  @{addr='AY4'; val='Trafikverkets inventeringar av arter i statliga vägmiljöer'}
)
foreach ($u in $textUpdates) {
    $c = $ws.Range($u.addr)
    $c.NumberFormat = "@"
    $c.Value = $u.val
    $c.ClearFormats()
}

# --- Cells that must end up blank in the new layout ---
$clearAddrs = @(
  "K2", "L2", "M2", "AH2", "AY2",
  "N3", "AF3",
  "I4", "J4", "K4", "L4", "M4", "N4", "AF4"
)
foreach ($a in $clearAddrs) {
    $ws.Range($a).ClearContents()
}

Write-Output "Row rotation (2 <- 3 <- 4 <- 2) applied"
